# Apply updates to publications worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publications")

# 1. Row 2: add missing tags value.
$ws.Range("K2").Value = "smartGrid, phaseTransition"

# 2. Row 20: "Emergent chirality..." (Tan et al.) paper moves from preprint to published (PRX Life).
$ws.Range("A20").Value = 45512
$ws.Range("D20").Value = "PRX Life"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = "033006"
$ws.Range("H20").Value = "tan2024emergent"
$ws.Range("I20").Value = "https://journals.aps.org/prxlife/abstract/10.1103/PRXLife.2.033006"
$ws.Hyperlinks.Add($ws.Range("I20"), "https://journals.aps.org/prxlife/abstract/10.1103/PRXLife.2.033006") | Out-Null
$ws.Range("I20").Style = $ws.Range("I16").Style

# 3. Row 21: "Active shape programming..." (Fuhrmann et al.) paper moves from preprint to published (Sci. Adv.).
$ws.Range("A21").Value = 45513
$ws.Range("B21").Value = "Active shape programming drives Drosophila wing disc eversion "
$ws.Range("C21").Value = "J. F. Fuhrmann, A. Krishna, J. Paijmans, C. Duclut, G. Cwikla, S. Eaton, M. Popović, F. Jülicher, C. D. Modes, N. A. Dye"
$ws.Range("D21").Value = "Sci. Adv."
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = "32"
$ws.Range("H21").Value = "fuhrmann2024active"
$ws.Range("I21").Value = "https://www.science.org/doi/10.1126/sciadv.adp0860"

# 4. Update selected cell in sheet view.
$ws.Range("B24").Select()
